# Feat 77: add metadata_dir field to excel template and validated jobs table
#
# Inserts a new "metadata_dir" column before the existing "modality0" column
# (i.e. at column D), shifting modality0 / modality0.source / modality1 /
# modality1.source one column to the right (D->E, E->F, F->G, G->H), and
# populates the header + sample rows for the new column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at D; this shifts D:G -> E:H and widens the sheet to
# column H. Excel also grows the dataValidation sqref (D2:D20/F2:F20 ->
# E2:E20/G2:G20) and the used-range dimension automatically.
$ws.Columns("D").Insert()

# Match the column formatting (width) used by its neighbours.
$ws.Columns("D").ColumnWidth = 12.17

# Header
$ws.Range("D1").Value = "metadata_dir"

# Sample data
$ws.Range("D2").Value = "/allen/aind/stage/fake/metadata_dir"
$ws.Range("D3").Value = "/allen/aind/stage/fake/Config"
# Row 4 intentionally left blank in the metadata_dir column.
